$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs", "Efna5", "Epha4", "ECs", 1, 0.3333333333333333, 0.09353433333333333, 0.280603, 0.0471360296668244, 0.0471360296668244, 3, 1, 7.189414999999999, 21.568245, 0.3327435204651371, 0.3327435204651371, 0.6724571390816666, 6.052114251735, 0.01568420845208829, 0.01568420845208829),
  @("ECs", "Efna5", "Epha4", "FAPs", 1, 0.3333333333333333, 0.09353433333333333, 0.280603, 0.0471360296668244, 0.0471360296668244, 3, 1, 12.90571233333333, 38.717137, 0.5973075912162077, 0.5973075912162077, 1.207127199290111, 10.864144793611, 0.02815470833978659, 0.02815470833978659),
  @("ECs", "Efna5", "Epha4", "M2", 1, 0.3333333333333333, 0.09353433333333333, 0.280603, 0.0471360296668244, 0.0471360296668244, 1, 0.3333333333333333, 0.02123933333333333, 0.063718, 0.0009830077336842938, 0.0009830077336842938, 0.001986606883777778, 0.017879461954, 0.0000463350816976607, 0.0000463350816976607),
  @("ECs", "Efna5", "Epha4", "sCs", 1, 0.3333333333333333, 0.09353433333333333, 0.280603, 0.0471360296668244, 0.0471360296668244, 3, 1, 1.490109666666666, 4.470329, 0.0689658805849709, 0.0689658805849709, 0.1393764142652222, 1.254387728387, 0.003250777793251857, 0.003250777793251857),
  @("FAPs", "Efna5", "Epha4", "ECs", 3, 1, 1.367717666666667, 4.103153, 0.6892525793933763, 0.6892525793933763, 3, 1, 7.189414999999999, 21.568245, 0.3327435204651371, 0.3327435204651371, 9.833089908498332, 88.49780917648499, 0.2293443297570284, 0.2293443297570284),
  @("FAPs", "Efna5", "Epha4", "FAPs", 3, 1, 1.367717666666667, 4.103153, 0.6892525793933763, 0.6892525793933763, 3, 1, 12.90571233333333, 38.717137, 0.5973075912162077, 0.5973075912162077, 17.65137075921789, 158.862336832961, 0.4116957979370155, 0.4116957979370155),
  @("FAPs", "Efna5", "Epha4", "M2", 3, 1, 1.367717666666667, 4.103153, 0.6892525793933763, 0.6892525793933763, 1, 0.3333333333333333, 0.02123933333333333, 0.063718, 0.0009830077336842938, 0.0009830077336842938, 0.02904941142822222, 0.261444702854, 0.0006775406160055365, 0.0006775406160055365),
  @("FAPs", "Efna5", "Epha4", "sCs", 3, 1, 1.367717666666667, 4.103153, 0.6892525793933763, 0.6892525793933763, 3, 1, 1.490109666666666, 4.470329, 0.0689658805849709, 0.0689658805849709, 2.038049316370778, 18.342443847337, 0.04753491108332676, 0.04753491108332676),
  @("sCs", "Efna5", "Epha4", "ECs", 3, 1, 0.523097, 1.569291, 0.2636113909397994, 0.2636113909397994, 3, 1, 7.189414999999999, 21.568245, 0.3327435204651371, 0.3327435204651371, 3.760761418255, 33.84685276429499, 0.08771498225602038, 0.08771498225602038),
  @("sCs", "Efna5", "Epha4", "FAPs", 3, 1, 0.523097, 1.569291, 0.2636113909397994, 0.2636113909397994, 3, 1, 12.90571233333333, 38.717137, 0.5973075912162077, 0.5973075912162077, 6.750939404429667, 60.758454639867, 0.1574570849394056, 0.1574570849394056),
  @("sCs", "Efna5", "Epha4", "M2", 3, 1, 0.523097, 1.569291, 0.2636113909397994, 0.2636113909397994, 1, 0.3333333333333333, 0.02123933333333333, 0.063718, 0.0009830077336842938, 0.0009830077336842938, 0.01111023154866667, 0.09999208393799999, 0.0002591320359810966, 0.0002591320359810966),
  @("sCs", "Efna5", "Epha4", "sCs", 3, 1, 0.523097, 1.569291, 0.2636113909397994, 0.2636113909397994, 3, 1, 1.490109666666666, 4.470329, 0.0689658805849709, 0.0689658805849709, 0.7794718963043332, 7.015247066739, 0.01818019170839228, 0.01818019170839228)
)

$r = 2
foreach ($row in $rows) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}
